# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the handoff
# package is ready: status text flips from "In Translation" to
# "Ready for handoff" and the associated generation timestamps are
# bumped forward. The status/date columns are also widened slightly
# to comfortably fit the new, longer status string.

$wb = $excel.ActiveWorkbook

$sheetOverview = $wb.Worksheets.Item("Overview")
$sheetZhCn     = $wb.Worksheets.Item("zh-cn")
$sheetDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -------------
$sheetOverview.Range("E2").Value = "Ready for handoff"
$sheetOverview.Range("F2").Value = "Ready for handoff"
$sheetZhCn.Range("C2").Value     = "Ready for handoff"
$sheetDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps ---------------------------------------------------------
# Overview!G2 and de-de!H2 shared "2016-09-06 17:38:57" -> bump to 17:39:47
$sheetOverview.Range("G2").Value = "2016-09-06 17:39:47"
$sheetDeDe.Range("H2").Value     = "2016-09-06 17:39:47"

# zh-cn!H2 "2016-09-06 17:38:52" -> bump to 17:39:42
$sheetZhCn.Range("H2").Value = "2016-09-06 17:39:42"

# --- Widen the Status columns (the grid quantizes ColumnWidth to the
# nearest 1/6th of a character, so we pick the input that lands on the
# closest achievable grid point to the target stored width of
# ~17.216 -> 17.1666...) --------------------------------------------
$sheetOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$sheetOverview.Columns.Item(6).ColumnWidth = 16.333333333333332
$sheetZhCn.Columns.Item(3).ColumnWidth     = 16.333333333333332
$sheetDeDe.Columns.Item(3).ColumnWidth     = 16.333333333333332
